{"js": "// Applies the LOM3246 doc update:\n//   1. Bump the \"Ativa\u00e7\u00e3o\" date from 01/01/2012 to 01/01/2023.\n//   2. Add an italic English translation paragraph after the Portuguese\n//      \"Objetivos\" paragraph.\n//   3. Add an italic English translation paragraph after the Portuguese\n//      \"Programa resumido\" paragraph.\n//   4. Collapse the four <w:br/>-separated runs in the \"Programa\" paragraph\n//      into a single run/text, then add an italic English translation\n//      paragraph right after it.\n\nconst body = context.document.body;\n\n// 1. Ativa\u00e7\u00e3o date.\nconst dateHits = body.search(\"Ativa\u00e7\u00e3o: 01/01/2012\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\ndateHits.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2023\", \"Replace\");\nawait context.sync();\n\n// 2. Objetivos -> add italic EN paragraph right after it.\nconst objHits = body.search(\n  \"Fornecer ao aluno o conhecimento das principais t\u00e9cnicas de caracteriza\u00e7\u00e3o f\u00edsica e qu\u00edmica de materiais.\",\n  { matchCase: true }\n);\nobjHits.load(\"items\");\nawait context.sync();\nconst objPara = objHits.items[0].paragraphs.getFirst();\nconst objEnPara = objPara.insertParagraph(\n  \"Provide the student with knowledge of the main techniques of physical and chemical characterization of materials.\",\n  \"After\"\n);\nobjEnPara.font.italic = true;\nawait context.sync();\n\n// 3. Programa resumido -> add italic EN paragraph right after it.\nconst summaryHits = body.search(\n  \"An\u00e1lise granulom\u00e9trica e superficial. An\u00e1lises microestruturais. An\u00e1lises t\u00e9rmicas. Reometria.\",\n  { matchCase: true }\n);\nsummaryHits.load(\"items\");\nawait context.sync();\nconst summaryPara = summaryHits.items[0].paragraphs.getFirst();\nconst summaryEnPara = summaryPara.insertParagraph(\n  \"Granulometric and surface analysis. Microstructural analyses. Thermal analysis. Rheometry.\",\n  \"After\"\n);\nsummaryEnPara.font.italic = true;\nawait context.sync();\n\n// 4. Programa -> merge the 4 <w:br/>-joined runs into one run of text,\n//    then add the italic EN paragraph right after it.\nconst programaHits = body.search(\n  \"An\u00e1lise granulom\u00e9trica. Adsor\u00e7\u00e3o BET, porosidade e picnometria.\",\n  { matchCase: true }\n);\nprogramaHits.load(\"items\");\nawait context.sync();\nconst programaPara = programaHits.items[0].paragraphs.getFirst();\n\nconst mergedPt =\n  \"An\u00e1lise granulom\u00e9trica. Adsor\u00e7\u00e3o BET, porosidade e picnometria. \" +\n  \"An\u00e1lises microestruturais: difra\u00e7\u00e3o de raios X, figura de Laue; \" +\n  \"espalhamento de raios X (SAXS). Difra\u00e7\u00e3o de el\u00e9trons. Microscopia \" +\n  \"\u00d3ptica. Microscopia eletr\u00f4nica, microan\u00e1lise de raios X (EDX e WDX). \" +\n  \"An\u00e1lises t\u00e9rmicas: An\u00e1lise t\u00e9rmica diferencial (DTA), calorimetria \" +\n  \"explorat\u00f3ria diferencial (DSC) e termogravimetria (TGA).\" +\n  \"Reometria de l\u00edquidos, solu\u00e7\u00f5es e pastas.\";\n\n// Clear the paragraph first so no stray xml:space=\"preserve\" survives from\n// the old (now gone) trailing-space run, then write the merged text back as\n// a single run.\nprogramaPara.clear();\nawait context.sync();\nprogramaPara.getRange().insertText(mergedPt, \"Replace\");\nawait context.sync();\n\nconst programaEn =\n  \"Grain size analysis. BET adsorption, porosity and pycnometry.\" +\n  \"Microstructural analysis: X-ray diffraction, Laue figure; X-ray \" +\n  \"scattering (SAXS). Electron diffraction. Optical Microscopy. Electron \" +\n  \"microscopy, X-ray microanalysis (EDX and WDX).Thermal analysis: \" +\n  \"Differential thermal analysis (DTA), differential scanning calorimetry \" +\n  \"(DSC) and thermogravimetry (TGA).Rheometry of liquids, solutions and \" +\n  \"pastes.\";\n\nconst programaEnPara = programaPara.insertParagraph(programaEn, \"After\");\nprogramaEnPara.font.italic = true;\nawait context.sync();\n", "ps1": "# Applies the LOM3246 doc update:\n#   1. Bump the \"Ativa\u00e7\u00e3o\" date from 01/01/2012 to 01/01/2023.\n#   2. Add an italic English translation paragraph after the Portuguese\n#      \"Objetivos\" paragraph.\n#   3. Add an italic English translation paragraph after the Portuguese\n#      \"Programa resumido\" paragraph.\n#   4. Collapse the four <w:br/>-separated runs in the \"Programa\" paragraph\n#      into a single run/text, then add an italic English translation\n#      paragraph right after it.\n\n$d = $word.ActiveDocument\n\n# 1. Ativa\u00e7\u00e3o date.\n$find = $d.Content.Find\n$find.Text = \"Ativa\u00e7\u00e3o: 01/01/2012\"\n$find.Replacement.Text = \"Ativa\u00e7\u00e3o: 01/01/2023\"\n[void]$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2. Objetivos -> add italic EN paragraph right after it.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Fornecer ao aluno o conhecimento das principais t\u00e9cnicas de caracteriza\u00e7\u00e3o f\u00edsica e qu\u00edmica de materiais.*\") {\n        $r = $p.Range\n        [void]$r.InsertParagraphAfter()\n        $enPara = $p.Next()\n        $enRange = $enPara.Range\n        $enRange.Text = \"Provide the student with knowledge of the main techniques of physical and chemical characterization of materials.\"\n        $enRange2 = $enPara.Range\n        [void]$enRange2.MoveEnd(1, -1)\n        $enRange2.Italic = 1\n        break\n    }\n}\n\n# 3. Programa resumido -> add italic EN paragraph right after it.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*An\u00e1lise granulom\u00e9trica e superficial. An\u00e1lises microestruturais. An\u00e1lises t\u00e9rmicas. Reometria.*\") {\n        $r = $p.Range\n        [void]$r.InsertParagraphAfter()\n        $enPara = $p.Next()\n        $enRange = $enPara.Range\n        $enRange.Text = \"Granulometric and surface analysis. Microstructural analyses. Thermal analysis. Rheometry.\"\n        $enRange2 = $enPara.Range\n        [void]$enRange2.MoveEnd(1, -1)\n        $enRange2.Italic = 1\n        break\n    }\n}\n\n# 4. Programa -> merge the 4 <w:br/>-joined runs into one run of text, then\n#    add the italic EN paragraph right after it.\n$mergedPt = \"An\u00e1lise granulom\u00e9trica. Adsor\u00e7\u00e3o BET, porosidade e picnometria. \" + `\n    \"An\u00e1lises microestruturais: difra\u00e7\u00e3o de raios X, figura de Laue; \" + `\n    \"espalhamento de raios X (SAXS). Difra\u00e7\u00e3o de el\u00e9trons. Microscopia \" + `\n    \"\u00d3ptica. Microscopia eletr\u00f4nica, microan\u00e1lise de raios X (EDX e WDX). \" + `\n    \"An\u00e1lises t\u00e9rmicas: An\u00e1lise t\u00e9rmica diferencial (DTA), calorimetria \" + `\n    \"explorat\u00f3ria diferencial (DSC) e termogravimetria (TGA).\" + `\n    \"Reometria de l\u00edquidos, solu\u00e7\u00f5es e pastas.\"\n\n$programaEn = \"Grain size analysis. BET adsorption, porosity and pycnometry.\" + `\n    \"Microstructural analysis: X-ray diffraction, Laue figure; X-ray \" + `\n    \"scattering (SAXS). Electron diffraction. Optical Microscopy. Electron \" + `\n    \"microscopy, X-ray microanalysis (EDX and WDX).Thermal analysis: \" + `\n    \"Differential thermal analysis (DTA), differential scanning calorimetry \" + `\n    \"(DSC) and thermogravimetry (TGA).Rheometry of liquids, solutions and \" + `\n    \"pastes.\"\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*An\u00e1lise granulom\u00e9trica. Adsor\u00e7\u00e3o BET, porosidade e picnometria*\") {\n        $r = $p.Range\n        [void]$r.MoveEnd(1, -1)\n        $r.Text = $mergedPt\n\n        $pr = $p.Range\n        [void]$pr.InsertParagraphAfter()\n        $enPara = $p.Next()\n        $enRange = $enPara.Range\n        $enRange.Text = $programaEn\n        $enRange2 = $enPara.Range\n        [void]$enRange2.MoveEnd(1, -1)\n        $enRange2.Italic = 1\n        break\n    }\n}\n"}
